$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sprint-Backlog")
$ws2 = $wb.Worksheets.Item("Burndown-Chart")

# --- Sprint-Backlog data updates (raw values; P5/Q5/R5 are formulas and recalc automatically) ---

$ws1.Range("P7").Value = 3
$ws1.Range("Q7").Value = 2
$ws1.Range("R7").Value = 2

$ws1.Range("Q8").Value = 0
$ws1.Range("R8").Value = 0

$ws1.Range("P11").Value = 3
$ws1.Range("Q11").Value = 5
$ws1.Range("R11").Value = 0

$ws1.Range("Q12").Value = 0
$ws1.Range("R12").Value = 0

$ws1.Range("Q13").Value = 0
$ws1.Range("R13").Value = 0

$ws1.Range("Q15").Value = 4
$ws1.Range("R15").Value = 3

$ws1.Range("Q16").Value = 0
$ws1.Range("R16").Value = 0

$ws1.Range("M31").Value = 3
$ws1.Range("N31").Value = 5
$ws1.Range("O31").Value = 0
$ws1.Range("P31").Value = 0
$ws1.Range("Q31").Value = 0
$ws1.Range("R31").Value = 0

$ws1.Range("N32").Value = 0
$ws1.Range("O32").Value = 0
$ws1.Range("P32").Value = 0
$ws1.Range("Q32").Value = 0
$ws1.Range("R32").Value = 0

$ws1.Range("O33").Value = 0
$ws1.Range("P33").Value = 0
$ws1.Range("Q33").Value = 0
$ws1.Range("R33").Value = 0

$ws1.Range("N35").Value = 0
$ws1.Range("O35").Value = 0
$ws1.Range("P35").Value = 0
$ws1.Range("Q35").Value = 0
$ws1.Range("R35").Value = 0

$ws1.Range("N36").Value = 0
$ws1.Range("O36").Value = 0
$ws1.Range("P36").Value = 0
$ws1.Range("Q36").Value = 0
$ws1.Range("R36").Value = 0

$ws1.Range("N37").Value = 0
$ws1.Range("O37").Value = 0
$ws1.Range("P37").Value = 0
$ws1.Range("Q37").Value = 0
$ws1.Range("R37").Value = 0

$ws1.Range("P67").Value = 4
$ws1.Range("Q67").Value = 4
$ws1.Range("R67").Value = 0

$ws1.Range("P68").Value = 0
$ws1.Range("Q68").Value = 0
$ws1.Range("R68").Value = 0

$ws1.Range("P69").Value = 0
$ws1.Range("Q69").Value = 0
$ws1.Range("R69").Value = 0

$ws1.Range("P71").Value = 4
$ws1.Range("Q71").Value = 4
$ws1.Range("R71").Value = 0

$ws1.Range("Q72").Value = 0
$ws1.Range("R72").Value = 0

$ws1.Range("Q73").Value = 0
$ws1.Range("R73").Value = 0

$ws1.Range("P75").Value = 4
$ws1.Range("Q75").Value = 4
$ws1.Range("R75").Value = 0

$ws1.Range("Q76").Value = 0
$ws1.Range("R76").Value = 0

$ws1.Range("Q77").Value = 0
$ws1.Range("R77").Value = 0

$ws1.Range("P87").Value = 4
$ws1.Range("Q87").Value = 2
$ws1.Range("R87").Value = 2

$ws1.Range("Q88").Value = 0
$ws1.Range("R88").Value = 0

$ws1.Range("Q91").Value = 6
$ws1.Range("R91").Value = 2

$ws1.Range("Q92").Value = 0
$ws1.Range("R92").Value = 0

$ws1.Range("Q97").Value = 4

$ws1.Range("P103").Value = 4
$ws1.Range("Q103").Value = 3
$ws1.Range("R103").Value = 1

$ws1.Range("Q104").Value = 0
$ws1.Range("R104").Value = 0

# --- Selection / view updates ---
# Touch Burndown-Chart's selection/scroll first (resets its stored topLeftCell),
# then finish on Sprint-Backlog so it remains the active/tab-selected sheet.
$ws2.Range("C23").Select()
$ws1.Range("P43").Select()

$excel.CalculateFull()
